# Add a new "2022-Q1" sheet (holdings detail) right before the "总计"
# (total) sheet, and prepend a matching summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new worksheet just before "总计" and rename it.
#    NOTE: worksheet variables in this host resolve positionally, so the
#    "总计" handle grabbed before the insert would silently start
#    pointing at the new sheet afterwards. Re-fetch it by name once the
#    insert/rename is done.
# ---------------------------------------------------------------------
$totalSheetBefore = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheetBefore)
$q1.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")

# Use an existing quarter sheet as a formatting template (header style,
# index-column style) so the new sheet matches the look of its peers.
$template = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Holdings detail rows.
# ---------------------------------------------------------------------
$rows = @(
    @{ idx = 0; code = "011056"; name = "博时汇兴回报一年持有期灵活配置混合"; size = "107.57"; pos = "67.02"; pct = "2.21"; mv = "2.3773"; rank = 10 },
    @{ idx = 1; code = "610004"; name = "信达澳银中小盘混合";                 size = "9.71";   pos = "93.54"; pct = "3.65"; mv = "0.3544"; rank = 10 },
    @{ idx = 2; code = "003877"; name = "富国久利稳健配置混合A";              size = "0.24";   pos = "26.60"; pct = "1.55"; mv = "0.0037"; rank = 8 },
    @{ idx = 3; code = "003878"; name = "富国久利稳健配置混合C";              size = "0.06";   pos = "26.60"; pct = "1.55"; mv = "0.0009"; rank = 8 }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q1.Cells.Item($r, 1).Value = $data.idx

    $q1.Cells.Item($r, 2).Value = "'" + $data.code
    $q1.Cells.Item($r, 2).ClearFormats()

    $q1.Cells.Item($r, 3).Value = $data.name

    $q1.Cells.Item($r, 4).Value = "'" + $data.size
    $q1.Cells.Item($r, 4).ClearFormats()

    $q1.Cells.Item($r, 5).Value = "'" + $data.pos
    $q1.Cells.Item($r, 5).ClearFormats()

    $q1.Cells.Item($r, 6).Value = "'" + $data.pct
    $q1.Cells.Item($r, 6).ClearFormats()

    $q1.Cells.Item($r, 7).Value = "'" + $data.mv
    $q1.Cells.Item($r, 7).ClearFormats()

    $q1.Cells.Item($r, 8).Value = $data.rank
}

# Re-apply the index-column (A) style used on every other quarter sheet.
$template.Range("A2").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. Prepend a "2022-Q1" summary row to the "总计" sheet.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 4
$totalSheet.Cells.Item(2, 4).Value = 2.74

$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)

# Renumber the 0-based index column for the rows that shifted down.
$lastRow = $totalSheet.UsedRange.Rows.Count
for ($r = 3; $r -le $lastRow; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
